$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as text, avoiding Excel auto-converting numeric-looking strings
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "43.714.14"
Set-TextCell $ws.Range("E2") "  +2.00%  "

# Row 3
Set-TextCell $ws.Range("D3") "2.211.51"
Set-TextCell $ws.Range("E3") "  -0.01%  "

# Row 4
Set-TextCell $ws.Range("E4") "  +0.04%  "

# Row 5
Set-TextCell $ws.Range("D5") "264.16"
Set-TextCell $ws.Range("E5") "  +2.75%  "

# Row 6
Set-TextCell $ws.Range("E6") "  +12.34%  "

# Row 7
Set-TextCell $ws.Range("E7") "  +0.86%  "

# Row 8
Set-TextCell $ws.Range("E8") "  -0.05%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.605"
Set-TextCell $ws.Range("E9") "  +1.71%  "

# Row 10
Set-TextCell $ws.Range("D10") "46.28"
Set-TextCell $ws.Range("E10") "  +9.42%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.0919"
Set-TextCell $ws.Range("E11") "  +1.24%  "

# Row 12
Set-TextCell $ws.Range("E12") "  +8.79%  "

# Row 13
Set-TextCell $ws.Range("E13") "  +2.05%  "

# Row 14
Set-TextCell $ws.Range("D14") "2.545.30"
Set-TextCell $ws.Range("E14") "  +0.14%  "

# Row 15
Set-TextCell $ws.Range("D15") "14.61"
Set-TextCell $ws.Range("E15") "  +0.77%  "

# Row 16
Set-TextCell $ws.Range("D16") "2.188.65"
Set-TextCell $ws.Range("E16") "  -1.47%  "

# Row 17
Set-TextCell $ws.Range("D17") "0.779"
Set-TextCell $ws.Range("E17") "  -0.62%  "

# Row 18
Set-TextCell $ws.Range("D18") "43.683.14"
Set-TextCell $ws.Range("E18") "  +1.98%  "

# Row 19
Set-TextCell $ws.Range("E19") "  +0.70%  "

# Row 20
Set-TextCell $ws.Range("D20") "5.97"
Set-TextCell $ws.Range("E20") "  +0.01%  "

# Row 21
Set-TextCell $ws.Range("D21") "70.05"
Set-TextCell $ws.Range("E21") "  -1.60%  "

# Row 22
Set-TextCell $ws.Range("E22") "  +7.84%  "

# Row 23
Set-TextCell $ws.Range("D23") "231.91"
Set-TextCell $ws.Range("E23") "  +0.74%  "

# Row 24
Set-TextCell $ws.Range("D24") "8.90"
Set-TextCell $ws.Range("E24") "  -5.36%  "

# Row 25
Set-TextCell $ws.Range("E25") "  +0.00%  "

# Row 26
Set-TextCell $ws.Range("B26") "PancakeSwap"
Set-TextCell $ws.Range("C26") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Range("D26") "2.43"
Set-TextCell $ws.Range("E26") "  +10.50%  "

# Row 27
Set-TextCell $ws.Range("B27") "Cosmos"
Set-TextCell $ws.Range("C27") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws.Range("D27") "10.83"
Set-TextCell $ws.Range("E27") "  +0.73%  "

# Row 28
Set-TextCell $ws.Range("B28") "WEMIXToken"
Set-TextCell $ws.Range("C28") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D28") "3.52"
Set-TextCell $ws.Range("E28") "  +5.58%  "

# Row 29
Set-TextCell $ws.Range("D29") "39.32"
Set-TextCell $ws.Range("E29") "  -8.74%  "

# Row 30
Set-TextCell $ws.Range("E30") "  +2.09%  "

# Row 31
Set-TextCell $ws.Range("D31") "175.10"
Set-TextCell $ws.Range("E31") "  +1.09%  "

# Row 32
Set-TextCell $ws.Range("E32") "  +2.54%  "

# Row 33
Set-TextCell $ws.Range("D33") "20.52"
Set-TextCell $ws.Range("E33") "  +0.86%  "

# Row 34
Set-TextCell $ws.Range("D34") "5.40"
Set-TextCell $ws.Range("E34") "  +3.32%  "

# Row 35
Set-TextCell $ws.Range("E35") "  +1.63%  "

# Row 36
Set-TextCell $ws.Range("D36") "0.110"
Set-TextCell $ws.Range("E36") "  +2.38%  "

# Row 37
Set-TextCell $ws.Range("D37") "0.0357"
Set-TextCell $ws.Range("E37") "  -1.27%  "

# Row 38
Set-TextCell $ws.Range("D38") "4.43"
Set-TextCell $ws.Range("E38") "  +1.83%  "

# Row 39
Set-TextCell $ws.Range("D39") "3.26"
Set-TextCell $ws.Range("E39") "  +15.64%  "

# Row 40
Set-TextCell $ws.Range("D40") "12.39"
Set-TextCell $ws.Range("E40") "  -4.18%  "

# Row 41
Set-TextCell $ws.Range("D41") "64.67"
Set-TextCell $ws.Range("E41") "  +7.65%  "

# Row 42
Set-TextCell $ws.Range("E42") "  -0.98%  "

# Row 43
Set-TextCell $ws.Range("D43") "5.54"
Set-TextCell $ws.Range("E43") "  +4.38%  "

# Row 44
Set-TextCell $ws.Range("D44") "0.204"
Set-TextCell $ws.Range("E44") "  +1.48%  "

# Row 45
Set-TextCell $ws.Range("D45") "100.59"
Set-TextCell $ws.Range("E45") "  -2.18%  "

# Row 46
Set-TextCell $ws.Range("E46") "  +0.58%  "

# Row 47
Set-TextCell $ws.Range("D47") "8.34"
Set-TextCell $ws.Range("E47") "  -0.14%  "

# Row 48
Set-TextCell $ws.Range("E48") "  +3.97%  "

# Row 49
Set-TextCell $ws.Range("B49") "ARBITRUM"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D49") "1.12"
Set-TextCell $ws.Range("E49") "  +0.60%  "

# Row 50
Set-TextCell $ws.Range("B50") "WOONetwork"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextCell $ws.Range("D50") "0.445"
Set-TextCell $ws.Range("E50") "  -3.97%  "

# Row 51
Set-TextCell $ws.Range("B51") "Stacks"
Set-TextCell $ws.Range("C51") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D51") "1.53"
Set-TextCell $ws.Range("E51") "  +6.93%  "
